# Adapt column header formatting to respective input file names:
#   "<Name>_old" -> "<Name>_FV2310"
#   "<Name>_new" -> "<Name>_FV2404"
# Wrap the sheet's data in an Excel Table, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the 20 header cells (columns A-J use the "old" format-version
#    suffix, columns L-U use the "new" one; column K just holds "diff").
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($baseNames[$i])_FV2310"
}
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($baseNames[$i])_FV2404"
}

# 2. Turn the used range A1:U70 into a proper Excel Table ("Table1"),
#    using the (now renamed) first row as the header row.
$range = $ws.Range("A1:U70")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# 3. Freeze panes at row 1 (so the header row stays visible while scrolling).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "Header renaming, table creation, and freeze panes applied."
